$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.807599666666667
$ws.Range("H2").Value = 5.422799
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 24.42119933333333
$ws.Range("N2").Value = 73.263598
$ws.Range("O2").Value = 0.4086816635579248
$ws.Range("P2").Value = 0.4086816635579248
$ws.Range("Q2").Value = 44.14375177453356
$ws.Range("R2").Value = 397.293765970802
$ws.Range("S2").Value = 0.4086816635579248
$ws.Range("T2").Value = 0.4086816635579248

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.807599666666667
$ws.Range("H3").Value = 5.422799
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 33.48129
$ws.Range("N3").Value = 100.44387
$ws.Range("O3").Value = 0.5602996441124273
$ws.Range("P3").Value = 0.5602996441124273
$ws.Range("Q3").Value = 60.52076864357
$ws.Range("R3").Value = 544.6869177921301
$ws.Range("S3").Value = 0.5602996441124273
$ws.Range("T3").Value = 0.5602996441124273

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.807599666666667
$ws.Range("H4").Value = 5.422799
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.853554333333333
$ws.Range("N4").Value = 5.560663
$ws.Range("O4").Value = 0.03101869232964781
$ws.Range("P4").Value = 0.03101869232964781
$ws.Range("Q4").Value = 3.350484195081889
$ws.Range("R4").Value = 30.154357755737
$ws.Range("S4").Value = 0.03101869232964781
$ws.Range("T4").Value = 0.03101869232964781
